$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3376
$ws.Range("F4").Value = 2476
$ws.Range("F5").Value = 337
$ws.Range("F7").Value = 1404
$ws.Range("F8").Value = 1104
$ws.Range("F9").Value = 310
$ws.Range("F10").Value = 521
$ws.Range("F12").Value = 16
$ws.Range("F15").Value = 8677
$ws.Range("F16").Value = 376
$ws.Range("F17").Value = 2489
$ws.Range("F18").Value = 260
$ws.Range("F22").Value = 595
$ws.Range("F24").Value = 1157
$ws.Range("F25").Value = 997
$ws.Range("F26").Value = 2032
$ws.Range("F27").Value = 2077
$ws.Range("F29").Value = 1767
$ws.Range("F33").Value = 121
$ws.Range("F34").Value = 53
$ws.Range("F35").Value = 97
$ws.Range("F36").Value = 187
$ws.Range("F37").Value = 7
$ws.Range("F38").Value = 307
$ws.Range("F40").Value = 250
$ws.Range("F41").Value = 434
$ws.Range("F42").Value = 708
$ws.Range("F43").Value = 78
$ws.Range("F44").Value = 265

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 18
$ws.Range("F4").Value = 7

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3376
$ws.Range("F4").Value = 2476
$ws.Range("F5").Value = 337
$ws.Range("F7").Value = 1404
$ws.Range("F9").Value = 1104
$ws.Range("F10").Value = 310
$ws.Range("F11").Value = 521
$ws.Range("F15").Value = 8677
$ws.Range("F16").Value = 376
$ws.Range("F17").Value = 2489
$ws.Range("F18").Value = 18
$ws.Range("F19").Value = 260
$ws.Range("F23").Value = 595
$ws.Range("F25").Value = 1157
$ws.Range("F26").Value = 997
$ws.Range("F27").Value = 2032
$ws.Range("F28").Value = 2077
$ws.Range("F29").Value = 1767
$ws.Range("F33").Value = 121
$ws.Range("F34").Value = 53
$ws.Range("F35").Value = 97
$ws.Range("F36").Value = 187
$ws.Range("F37").Value = 7
$ws.Range("F38").Value = 307
$ws.Range("F40").Value = 250
$ws.Range("F41").Value = 434
$ws.Range("F42").Value = 7
$ws.Range("F46").Value = 709
$ws.Range("F48").Value = 78
$ws.Range("F49").Value = 265
